$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2764.7778
$ws.Range("J32").Value = 3379
$ws.Range("L32").Value = 3379
$ws.Range("N32").Value = -4031
$ws.Range("H43").Value = 4252.2173
$ws.Range("I43").Value = 1800
$ws.Range("K43").Value = 1800
$ws.Range("M43").Value = -1731
$ws.Range("H88").Value = 26366322
$ws.Range("J88").Value = 3253751
$ws.Range("L88").Value = 3253751
$ws.Range("N88").Value = -3254563
$ws.Range("H91").Value = 26366322
$ws.Range("J91").Value = 3253751
$ws.Range("L91").Value = 3253751
$ws.Range("N91").Value = -3256559
$ws.Range("H106").Value = 3450
$ws.Range("J106").Value = 3450
$ws.Range("L106").Value = 3450
$ws.Range("N106").Value = -4712
$ws.Range("H112").Value = 168916
$ws.Range("I112").Value = 167666.17
$ws.Range("J112").Value = 170165.83
$ws.Range("K112").Value = 502998.51
$ws.Range("L112").Value = 510497.49
$ws.Range("M112").Value = -501890.51
$ws.Range("N112").Value = -512713.49
$ws.Range("H116").Value = 4538.0835
$ws.Range("I116").Value = 4382.5
$ws.Range("J116").Value = 6249.5
$ws.Range("K116").Value = 4382.5
$ws.Range("L116").Value = 6249.5
$ws.Range("M116").Value = -940.5
$ws.Range("N116").Value = -13133.5
$ws.Range("H135").Value = 23809978
$ws.Range("I135").Value = 26316116
$ws.Range("J135").Value = 1664.5
$ws.Range("K135").Value = 236845044
$ws.Range("L135").Value = 14980.5
$ws.Range("M135").Value = -236842509
$ws.Range("N135").Value = -20050.5
$ws.Range("H138").Value = 3399.647
$ws.Range("I138").Value = 2751
$ws.Range("K138").Value = 8253
$ws.Range("M138").Value = -3113

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1328.1818
$ws.Range("I2").Value = 576.1667
$ws.Range("K2").Value = 576.1667
$ws.Range("M2").Value = -463.1667
$ws.Range("H5").Value = 2719.1667
$ws.Range("I5").Value = 554.2857
$ws.Range("J5").Value = 4096.8184
$ws.Range("K5").Value = 554.2857
$ws.Range("L5").Value = 4096.8184
$ws.Range("M5").Value = -442.2857
$ws.Range("N5").Value = -4320.8184
$ws.Range("H32").Value = 3143.4082
$ws.Range("I32").Value = 1591.5682
$ws.Range("K32").Value = 1591.5682
$ws.Range("M32").Value = -1304.5682
$ws.Range("H45").Value = 4563.643
$ws.Range("I45").Value = 5559.2
$ws.Range("K45").Value = 5559.2
$ws.Range("M45").Value = -5182.2
$ws.Range("H61").Value = 45462320
$ws.Range("I61").Value = 62507396
$ws.Range("J61").Value = 8781.833000000001
$ws.Range("K61").Value = 62507396
$ws.Range("L61").Value = 8781.833000000001
$ws.Range("M61").Value = -62507184
$ws.Range("N61").Value = -9205.833000000001
$ws.Range("H98").Value = 40070.6
$ws.Range("J98").Value = 40070.6
$ws.Range("L98").Value = 40070.6
$ws.Range("N98").Value = -46060.6
$ws.Range("H116").Value = 1328.1818
$ws.Range("I116").Value = 576.1667
$ws.Range("K116").Value = 576.1667
$ws.Range("M116").Value = 1717.8333
$ws.Range("H136").Value = 45462320
$ws.Range("I136").Value = 62507396
$ws.Range("J136").Value = 8781.833000000001
$ws.Range("K136").Value = 187522188
$ws.Range("L136").Value = 26345.499
$ws.Range("M136").Value = -187519638
$ws.Range("N136").Value = -31445.499

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1328.1818
$ws.Range("I3").Value = 576.1667
$ws.Range("K3").Value = 576.1667
$ws.Range("M3").Value = -462.1667
$ws.Range("H4").Value = 2719.1667
$ws.Range("I4").Value = 554.2857
$ws.Range("J4").Value = 4096.8184
$ws.Range("K4").Value = 554.2857
$ws.Range("L4").Value = 4096.8184
$ws.Range("M4").Value = -439.2857
$ws.Range("N4").Value = -4326.8184
$ws.Range("H86").Value = 2664.8667
$ws.Range("I86").Value = 2849.7
$ws.Range("J86").Value = 2295.2
$ws.Range("K86").Value = 2849.7
$ws.Range("L86").Value = 2295.2
$ws.Range("M86").Value = -1726.7
$ws.Range("N86").Value = -4541.2
$ws.Range("H89").Value = 2664.8667
$ws.Range("I89").Value = 2849.7
$ws.Range("J89").Value = 2295.2
$ws.Range("K89").Value = 14248.5
$ws.Range("L89").Value = 11476
$ws.Range("M89").Value = -8632.5
$ws.Range("N89").Value = -22708
$ws.Range("H94").Value = 1785.303
$ws.Range("I94").Value = 1842.8518
$ws.Range("J94").Value = 1526.3334
$ws.Range("K94").Value = 1842.8518
$ws.Range("L94").Value = 1526.3334
$ws.Range("M94").Value = -1391.8518
$ws.Range("N94").Value = -2428.3334
$ws.Range("H95").Value = 17289.8
$ws.Range("J95").Value = 17289.8
$ws.Range("L95").Value = 17289.8
$ws.Range("N95").Value = -22781.8
$ws.Range("H105").Value = 3614.125
$ws.Range("I105").Value = 3195
$ws.Range("K105").Value = 3195
$ws.Range("M105").Value = -1448
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1428903.9
$ws.Range("I6").Value = 1428903.9
$ws.Range("K6").Value = 1428903.9
$ws.Range("M6").Value = -1428790.9
$ws.Range("H7").Value = 6132.0586
$ws.Range("I7").Value = 10238
$ws.Range("K7").Value = 10238
$ws.Range("M7").Value = -10125
$ws.Range("H31").Value = 6204.143
$ws.Range("I31").Value = 9757.625
$ws.Range("J31").Value = 1466.1666
$ws.Range("K31").Value = 9757.625
$ws.Range("L31").Value = 1466.1666
$ws.Range("M31").Value = -9462.625
$ws.Range("N31").Value = -2056.1666
$ws.Range("H34").Value = 6204.143
$ws.Range("I34").Value = 9757.625
$ws.Range("J34").Value = 1466.1666
$ws.Range("K34").Value = 9757.625
$ws.Range("L34").Value = 1466.1666
$ws.Range("M34").Value = -9555.625
$ws.Range("N34").Value = -1870.1666
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H94").Value = 2214.2
$ws.Range("I94").Value = 2737
$ws.Range("K94").Value = 2737
$ws.Range("M94").Value = -2286
$ws.Range("H134").Value = 10002431
$ws.Range("I134").Value = 11907190
$ws.Range("J134").Value = 2448.25
$ws.Range("K134").Value = 35721570
$ws.Range("L134").Value = 7344.75
$ws.Range("M134").Value = -35719035
$ws.Range("N134").Value = -12414.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 112214.336
$ws.Range("I113").Value = 125824.5
$ws.Range("K113").Value = 125824.5
$ws.Range("M113").Value = -123654.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4192.25
$ws.Range("I7").Value = 3881.5
$ws.Range("K7").Value = 3881.5
$ws.Range("M7").Value = -3769.5
$ws.Range("H68").Value = 1200
$ws.Range("J68").Value = 1200
$ws.Range("L68").Value = 1200
$ws.Range("N68").Value = -2698
$ws.Range("H71").Value = 1200
$ws.Range("J71").Value = 1200
$ws.Range("L71").Value = 6000
$ws.Range("N71").Value = -13488
$ws.Range("H101").Value = 18785.834
$ws.Range("J101").Value = 18785.834
$ws.Range("L101").Value = 18785.834
$ws.Range("N101").Value = -25275.834
$ws.Range("H122").Value = 7998.5
$ws.Range("I122").Value = 7998
$ws.Range("K122").Value = 23994
$ws.Range("M122").Value = -21544
$ws.Range("H126").Value = 4192.25
$ws.Range("I126").Value = 3881.5
$ws.Range("K126").Value = 11644.5
$ws.Range("M126").Value = -9174.5
$ws.Range("H132").Value = 21823596
$ws.Range("I132").Value = 25268346
$ws.Range("K132").Value = 75805038
$ws.Range("M132").Value = -75802508

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1838.7407
$ws.Range("I126").Value = 1689.8096
$ws.Range("K126").Value = 5069.4288
$ws.Range("M126").Value = -2599.4288
$ws.Range("H132").Value = 35715770
$ws.Range("I132").Value = 38463028
$ws.Range("K132").Value = 115389084
$ws.Range("M132").Value = -115386554
$ws.Range("H136").Value = 20001040
$ws.Range("I136").Value = 26316718
$ws.Range("K136").Value = 78950154
$ws.Range("M136").Value = -78947604
